# Add a "Revisado" marker column (I) for groups that have been reviewed,
# mirroring the merged "NETWORK" column (H) layout for the grouped rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Merge the I-column cells for the same row blocks that are merged in H
# (H18:H21 / H22:H24), so the "Revisado" label spans the whole group.
$ws.Range("I18:I21").Merge()
$ws.Range("I22:I24").Merge()

# Set the "Revisado" text on the group leaders and on the two standalone rows.
$ws.Range("I18").Value = "Revisado"
$ws.Range("I22").Value = "Revisado"
$ws.Range("I28").Value = "Revisado"
$ws.Range("I30").Value = "Revisado"

# Center-align column I for the merged group rows (matches column H styling).
$ws.Range("I18").HorizontalAlignment = -4108
$ws.Range("I19").HorizontalAlignment = -4108
$ws.Range("I20").HorizontalAlignment = -4108
$ws.Range("I21").HorizontalAlignment = -4108
$ws.Range("I22").HorizontalAlignment = -4108
$ws.Range("I23").HorizontalAlignment = -4108
$ws.Range("I24").HorizontalAlignment = -4108

# Move the active selection, as left by the author after the edit.
$ws.Range("I13").Select()
